$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "97.081.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.38%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.729.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.34%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.63%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.25%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "661.96"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.82%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.428"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.17%  "

$ws.Range("E9").Value = "  -1.82%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.01%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.729.61"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.43%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000328"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +22.02%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "45.18"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.89%  "

$ws.Range("E14").Value = "  +1.45%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.93"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.54%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.431.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.43%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "96.771.87"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.31%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.92%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.725.81"
$ws.Range("D19").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.54%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.18%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.510"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "528.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.21%  "

$ws.Range("E24").Value = "  +0.67%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000227"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +11.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.92%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "109.42"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.32%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.195"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +15.99%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "13.74"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.62%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.933.48"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.47%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.08"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.26%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.05"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.21%  "

$ws.Range("E33").Value = "  +0.16%  "

$ws.Range("E34").Value = "  +3.58%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.87"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.50%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "32.99"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.33%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.43%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "650.25"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.598"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.63%  "

$ws.Range("E40").Value = "  -0.05%  "

$ws.Range("E42").Value = "  +4.62%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.60%  "

$ws.Range("E44").Value = "  +3.87%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.04"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.34%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.988"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.06%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.481"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +9.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0460"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.40%  "

$ws.Range("E49").Value = "  +3.47%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.78"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.45%  "

$ws.Range("E51").Value = "  -0.29%  "
